$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need the "@" (Text) number
# format applied first, otherwise Excel auto-converts the text into a numeric
# value (losing the original fixed-decimal formatting, e.g. "0.000008618"
# would turn into "8.618E-06").
$textCells = @("D5", "D6", "D7", "D8", "D10", "D11", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated crypto price (column D) and 1h volume change (column E) values
$ws.Range("D2").Value = '26.909.87'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '1.843.20'
$ws.Range("E3").Value = '  +1.44%  '
$ws.Range("E4").Value = '  +0.71%  '
$ws.Range("D5").Value = '308.99'
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("D7").Value = '0.4685'
$ws.Range("E7").Value = '  +3.41%  '
$ws.Range("D8").Value = '0.3661'
$ws.Range("E8").Value = '  +1.72%  '
$ws.Range("E9").Value = '  +0.69%  '
$ws.Range("D10").Value = '0.9262'
$ws.Range("E10").Value = '  +2.79%  '
$ws.Range("D11").Value = '19.56'
$ws.Range("E11").Value = '  +0.87%  '
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").Value = '1.815.36'
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").Value = '5.283'
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("D15").Value = '6.397'
$ws.Range("E15").Value = '  +0.92%  '
$ws.Range("D16").Value = '88.19'
$ws.Range("E16").Value = '  +3.26%  '
$ws.Range("D17").Value = '1.009'
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("D18").Value = '0.000008618'
$ws.Range("E18").Value = '  +0.65%  '
$ws.Range("D19").Value = '1.006'
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("D20").Value = '26.969.51'
$ws.Range("E20").Value = '  +1.36%  '
$ws.Range("D21").Value = '14.43'
$ws.Range("E22").Value = '  +0.93%  '
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("D24").Value = '1.938'
$ws.Range("E24").Value = '  -0.99%  '
$ws.Range("D25").Value = '152.28'
$ws.Range("E25").Value = '  +1.00%  '
$ws.Range("E26").Value = '  +2.21%  '
$ws.Range("D27").Value = '2.023'
$ws.Range("E27").Value = '  -1.87%  '
$ws.Range("D28").Value = '114.16'
$ws.Range("E28").Value = '  +1.38%  '
$ws.Range("D29").Value = '4.877'
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("D30").Value = '0.08854'
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("D31").Value = '3.219'
$ws.Range("E31").Value = '  +2.98%  '
$ws.Range("D32").Value = '1.175'
$ws.Range("E32").Value = '  +5.65%  '
$ws.Range("D33").Value = '0.7455'
$ws.Range("E33").Value = '  -0.98%  '
$ws.Range("D34").Value = '2.790'
$ws.Range("E34").Value = '  +1.35%  '
$ws.Range("E35").Value = '  +0.66%  '
$ws.Range("E36").Value = '  +0.96%  '
$ws.Range("D37").Value = '0.01939'
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D38").Value = '2.966'
$ws.Range("E38").Value = '  +2.05%  '
$ws.Range("D39").Value = '0.05184'
$ws.Range("E39").Value = '  +1.41%  '
$ws.Range("D40").Value = '0.5205'
$ws.Range("E40").Value = '  +1.91%  '
$ws.Range("D41").Value = '6.908'
$ws.Range("E41").Value = '  +2.06%  '
$ws.Range("E42").Value = '  +0.70%  '
$ws.Range("D43").Value = '8.121'
$ws.Range("E43").Value = '  +0.66%  '
$ws.Range("D44").Value = '10.52'
$ws.Range("E44").Value = '  +4.45%  '
$ws.Range("D45").Value = '0.4689'
$ws.Range("E45").Value = '  -1.08%  '
$ws.Range("D46").Value = '1.006'
$ws.Range("E46").Value = '  +0.50%  '
$ws.Range("D47").Value = '100.43'
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("D48").Value = '1.607'
$ws.Range("E48").Value = '  +1.62%  '
$ws.Range("D49").Value = '65.46'
$ws.Range("E49").Value = '  +2.30%  '
$ws.Range("D50").Value = '0.06038'
$ws.Range("E50").Value = '  +0.87%  '
$ws.Range("D51").Value = '0.8900'
$ws.Range("E51").Value = '  +4.88%  '
